# Auto-generated script applying scheduled market-data refresh values
# to the Leve profit sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("H100").Value = 1911.75
$ws.Range("J100").Value = 2999.6667
$ws.Range("L100").Value = 2999.6667
$ws.Range("N100").Value = -4081.6667
$ws.Range("H106").Value = 2442.889
$ws.Range("I106").Value = 2442.889
$ws.Range("K106").Value = 2442.889
$ws.Range("M106").Value = -1811.889
$ws.Range("H125").Value = 916
$ws.Range("H132").Value = 1301.0625
$ws.Range("I132").Value = 1214
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 3642
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = -1112
$ws.Range("N132").Value = -17060
$ws.Range("H141").Value = 4998.75
$ws.Range("J141").Value = 5666.6665
$ws.Range("L141").Value = 16999.9995
$ws.Range("N141").Value = -27359.9995
$ws.Range("N75").ClearContents()
$ws.Range("N78").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3271.0518
$ws.Range("I32").Value = 2038
$ws.Range("J32").Value = 7146.357
$ws.Range("K32").Value = 2038
$ws.Range("L32").Value = 7146.357
$ws.Range("M32").Value = -1751
$ws.Range("N32").Value = -7720.357
$ws.Range("H88").Value = 2729.7334
$ws.Range("J88").Value = 3749.8333
$ws.Range("L88").Value = 3749.8333
$ws.Range("N88").Value = -4561.8333
$ws.Range("H91").Value = 2729.7334
$ws.Range("J91").Value = 3749.8333
$ws.Range("L91").Value = 3749.8333
$ws.Range("N91").Value = -6557.8333
$ws.Range("H97").Value = 484.83334
$ws.Range("I97").Value = 491.22223
$ws.Range("K97").Value = 491.22223
$ws.Range("M97").Value = 4.777769999999975
$ws.Range("H98").Value = 18250
$ws.Range("J98").Value = 18250
$ws.Range("L98").Value = 18250
$ws.Range("N98").Value = -24240
$ws.Range("H104").Value = 38999.332
$ws.Range("J104").Value = 38999.332
$ws.Range("L104").Value = 38999.332
$ws.Range("N104").Value = -45987.332
$ws.Range("H132").Value = 2909.3572
$ws.Range("I132").Value = 1970.4445
$ws.Range("K132").Value = 5911.333500000001
$ws.Range("M132").Value = -3381.333500000001
$ws.Range("H137").Value = 12000
$ws.Range("I137").Value = 12000
$ws.Range("K137").Value = 12000
$ws.Range("M137").Value = -6900

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1614.5625
$ws.Range("I20").Value = 1510.1666
$ws.Range("J20").Value = 1927.75
$ws.Range("K20").Value = 1510.1666
$ws.Range("L20").Value = 1927.75
$ws.Range("M20").Value = -1263.1666
$ws.Range("N20").Value = -2421.75
$ws.Range("H36").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("H86").Value = 252524.88
$ws.Range("I86").Value = 2699.8333
$ws.Range("K86").Value = 2699.8333
$ws.Range("M86").Value = -1576.8333
$ws.Range("H89").Value = 252524.88
$ws.Range("I89").Value = 2699.8333
$ws.Range("K89").Value = 13499.1665
$ws.Range("M89").Value = -7883.166499999999
$ws.Range("H132").Value = 90000
$ws.Range("J132").Value = 90000
$ws.Range("L132").Value = 90000
$ws.Range("N132").Value = -100120
$ws.Range("H134").Value = 8478.963
$ws.Range("I134").Value = 8931.904
$ws.Range("J134").Value = 6893.6665
$ws.Range("K134").Value = 26795.712
$ws.Range("L134").Value = 20680.9995
$ws.Range("M134").Value = -24260.712
$ws.Range("N134").Value = -25750.9995
$ws.Range("N36").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2799
$ws.Range("I31").Value = 1102.1111
$ws.Range("J31").Value = 4980.7144
$ws.Range("K31").Value = 1102.1111
$ws.Range("L31").Value = 4980.7144
$ws.Range("M31").Value = -807.1111000000001
$ws.Range("N31").Value = -5570.7144
$ws.Range("H34").Value = 2799
$ws.Range("I34").Value = 1102.1111
$ws.Range("J34").Value = 4980.7144
$ws.Range("K34").Value = 1102.1111
$ws.Range("L34").Value = 4980.7144
$ws.Range("M34").Value = -900.1111000000001
$ws.Range("N34").Value = -5384.7144
$ws.Range("H93").Value = 4933
$ws.Range("I93").Value = 4933
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 4933
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = -3061
$ws.Range("H96").Value = 23750
$ws.Range("J96").Value = 23750
$ws.Range("L96").Value = 23750
$ws.Range("N96").Value = -29242
$ws.Range("H122").Value = 1540.3871
$ws.Range("I122").Value = 1333.8096
$ws.Range("J122").Value = 1974.2
$ws.Range("K122").Value = 4001.4288
$ws.Range("L122").Value = 5922.6
$ws.Range("M122").Value = -1551.4288
$ws.Range("N122").Value = -10822.6
$ws.Range("N93").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H96").Value = 450
$ws.Range("I96").Value = 450
$ws.Range("K96").Value = 1350
$ws.Range("M96").Value = 709
$ws.Range("H98").Value = 674.25
$ws.Range("J98").Value = 648.5
$ws.Range("L98").Value = 1945.5
$ws.Range("N98").Value = -4941.5
$ws.Range("H107").Value = 610.6
$ws.Range("J107").Value = 638.3570999999999
$ws.Range("L107").Value = 1915.0713
$ws.Range("N107").Value = -5755.0713
$ws.Range("H116").Value = 2526.3333
$ws.Range("J116").Value = 2526.3333
$ws.Range("L116").Value = 7578.999899999999
$ws.Range("N116").Value = -14462.9999
$ws.Range("H129").Value = 35390.24
$ws.Range("I129").Value = 760.8570999999999
$ws.Range("J129").Value = 52704.93
$ws.Range("K129").Value = 2282.5713
$ws.Range("L129").Value = 158114.79
$ws.Range("M129").Value = 2717.4287
$ws.Range("N129").Value = -168114.79
$ws.Range("H131").Value = 6956265.5
$ws.Range("I131").Value = 83333830
$ws.Range("J131").Value = 12850.091
$ws.Range("K131").Value = 250001490
$ws.Range("L131").Value = 38550.273
$ws.Range("M131").Value = -249996450
$ws.Range("N131").Value = -48630.273

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5253.4443
$ws.Range("I70").Value = 5877.4
$ws.Range("J70").Value = 4473.5
$ws.Range("K70").Value = 5877.4
$ws.Range("L70").Value = 4473.5
$ws.Range("M70").Value = -5607.4
$ws.Range("N70").Value = -5013.5
$ws.Range("H73").Value = 5253.4443
$ws.Range("I73").Value = 5877.4
$ws.Range("J73").Value = 4473.5
$ws.Range("K73").Value = 5877.4
$ws.Range("L73").Value = 4473.5
$ws.Range("M73").Value = -4941.4
$ws.Range("N73").Value = -6345.5
$ws.Range("H80").Value = 2503.8572
$ws.Range("I80").Value = 2406.75
$ws.Range("K80").Value = 2406.75
$ws.Range("M80").Value = -1408.75
$ws.Range("H83").Value = 2503.8572
$ws.Range("I83").Value = 2406.75
$ws.Range("K83").Value = 12033.75
$ws.Range("M83").Value = -7041.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4720.091
$ws.Range("I7").Value = 2325.8
$ws.Range("J7").Value = 6715.3335
$ws.Range("K7").Value = 2325.8
$ws.Range("L7").Value = 6715.3335
$ws.Range("M7").Value = -2213.8
$ws.Range("N7").Value = -6939.3335
$ws.Range("H40").Value = 8175.778
$ws.Range("I40").Value = 5352.857
$ws.Range("J40").Value = 9972.182000000001
$ws.Range("K40").Value = 5352.857
$ws.Range("L40").Value = 9972.182000000001
$ws.Range("M40").Value = -5216.857
$ws.Range("N40").Value = -10244.182
$ws.Range("H56").Value = 3649.6667
$ws.Range("I56").Value = 3974.5
$ws.Range("J56").Value = 3000
$ws.Range("K56").Value = 3974.5
$ws.Range("L56").Value = 3000
$ws.Range("M56").Value = -3283.5
$ws.Range("N56").Value = -4382
$ws.Range("H61").Value = 2717.1667
$ws.Range("I61").Value = 2260.6
$ws.Range("K61").Value = 2260.6
$ws.Range("M61").Value = -2058.6
$ws.Range("H82").Value = 3051.9333
$ws.Range("J82").Value = 3942.2222
$ws.Range("L82").Value = 3942.2222
$ws.Range("N82").Value = -4664.2222
$ws.Range("H85").Value = 3051.9333
$ws.Range("J85").Value = 3942.2222
$ws.Range("L85").Value = 3942.2222
$ws.Range("N85").Value = -6438.2222
$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("H113").Value = 2717.1667
$ws.Range("I113").Value = 2260.6
$ws.Range("K113").Value = 2260.6
$ws.Range("M113").Value = -90.59999999999991
$ws.Range("H126").Value = 4720.091
$ws.Range("I126").Value = 2325.8
$ws.Range("J126").Value = 6715.3335
$ws.Range("K126").Value = 6977.400000000001
$ws.Range("L126").Value = 20146.0005
$ws.Range("M126").Value = -4507.400000000001
$ws.Range("N126").Value = -25086.0005
$ws.Range("H132").Value = 2089.15
$ws.Range("J132").Value = 2164.3125
$ws.Range("L132").Value = 6492.9375
$ws.Range("N132").Value = -11552.9375
$ws.Range("N96").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H55").Value = 11000
$ws.Range("I55").Value = 2000
$ws.Range("K55").Value = 2000
$ws.Range("M55").Value = -1723
$ws.Range("H126").Value = 4882.9614
$ws.Range("I126").Value = 3870.111
$ws.Range("K126").Value = 11610.333
$ws.Range("M126").Value = -9140.332999999999
